$d = $word.ActiveDocument

# Locate the anchor paragraphs dynamically by their text content, so the
# script is resilient to exact paragraph numbering.
#
# We want to remove everything that follows the first product entry
# (the one ending in "Lugar: Frutal, Minas Gerais.") up to and including
# the second product entry's "Lugar: Birigui, Sao Paulo." line (the
# "GPSOM BIRIGUI" store entry), while leaving the trailing separator
# line and empty paragraph untouched.

$frutalPara = $null
$biriguiCount = 0
$secondBiriguiPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text

    if ($frutalPara -eq $null -and $txt -like "Lugar:*Frutal*") {
        $frutalPara = $p
    }

    if ($txt -like "Lugar:*Birigui*") {
        $biriguiCount = $biriguiCount + 1
        if ($biriguiCount -eq 2) {
            $secondBiriguiPara = $p
        }
    }
}

if ($frutalPara -ne $null -and $secondBiriguiPara -ne $null) {
    # Delete from right after "Lugar: Frutal, Minas Gerais." paragraph
    # through the end of the second "Lugar: Birigui, Sao Paulo." paragraph.
    $deleteRange = $d.Range($frutalPara.Range.End, $secondBiriguiPara.Range.End)
    $deleteRange.Delete()
}

$d.Save()
